# Weekly price-sheet update: a new, most-recent record is inserted at
# row 105 (pushing the existing historical rows 105-116 down to 106-117).
#
# Equivalent to Excel's Insert > Entire Row at row 105.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(105).Insert()

# Populate the newly-inserted row with this week's record.
$ws.Cells.Item(105, 1).Value  = 11
$ws.Cells.Item(105, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(105, 3).Value  = "Bíobío"
$ws.Cells.Item(105, 4).Value  = 44474
$ws.Cells.Item(105, 5).Value  = 8
$ws.Cells.Item(105, 6).Value  = 100114001
$ws.Cells.Item(105, 7).Value  = "Papa"
$ws.Cells.Item(105, 8).Value  = "Patagonia"
$ws.Cells.Item(105, 9).Value  = "1a (guarda)"
$ws.Cells.Item(105, 10).Value = 2000
$ws.Cells.Item(105, 11).Value = 11000
$ws.Cells.Item(105, 12).Value = 11500
$ws.Cells.Item(105, 13).Value = 11250
$ws.Cells.Item(105, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(105, 15).Value = "Provincia de Arauco"
$ws.Cells.Item(105, 16).Value = 450
$ws.Cells.Item(105, 17).Value = 25
$ws.Cells.Item(105, 18).Value = "Hortaliza"
